# xl data change commit
#
# - "customer service page title" row (B5) gets new display text: the
#   old shared string "Amazon Customer Service Support – Amazon.com" is
#   replaced by "Help & Contact Us - Amazon Customer Service".
# - "amazon home page url" row (B2) becomes a live hyperlink pointing at
#   its own URL text, mirroring the existing hyperlink already present
#   on the "customer service page url" row (B6).
# - Active selection ends up on B5 (the cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the customer-service title cell text.
$ws.Range("B5").Value = "Help & Contact Us - Amazon Customer Service"

# Turn B2 (the home page URL) into a hyperlink to itself, same pattern
# as the existing B6 hyperlink.
$ws.Hyperlinks.Add($ws.Range("B2"), $ws.Range("B2").Text)

# Leave the active cell/selection on B5.
$ws.Range("B5").Select()
